# stats.xlsx edit — 2021-07-27
# Summary of the change (per the OOXML diff):
#  - Sheet1!A7: "Simple" -> "All"
#  - Sheet1 becomes the active/selected tab; selection moves to B16
#  - correlation (Sheet2):
#      * a new title row is inserted above the existing "4184 words:" table
#        (shifts the existing 5-row block down by one row)
#      * a brand new second table ("all degree >= 5 (3406 words):") is
#        appended a few rows below, with the same column headers and new
#        numbers
#      * correlation is no longer the active/selected tab; selection moves
#        to E17

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- correlation sheet: push the existing table down and add its title ---
# Insert the blank row first (while the sheet is still in its original
# 5-row shape) so later writes below row 11 land on the right rows.
$ws2.Rows.Item(1).Insert()

# (Do the shared-string-producing writes in the same order the new
# strings appear in the saved workbook: "all degree ..." title, then
# "4184 words:" title, then Sheet1's "All".)

# --- correlation sheet: brand new second table beneath the first ---
$ws2.Range("A12").Value = "all degree >= 5 (3406 words):"

$ws2.Range("A1").Value = "4184 words:"

$ws2.Range("B13").Value = "all"
$ws2.Range("C13").Value = "replace"
$ws2.Range("D13").Value = "delete_insert"
$ws2.Range("E13").Value = "pos<50%"

$ws2.Range("A14").Value = "replace"
$ws2.Range("B14").Value = 0.96

$ws2.Range("A15").Value = "delete_insert"
$ws2.Range("B15").Value = 0.54
$ws2.Range("C15").Value = 0.28000000000000003

$ws2.Range("A16").Value = "pos<50%"
$ws2.Range("B16").Value = 0.92
$ws2.Range("C16").Value = 0.9
$ws2.Range("D16").Value = 0.46

$ws2.Range("A17").Value = "pos>=50%"
$ws2.Range("B17").Value = 0.72
$ws2.Range("C17").Value = 0.66
$ws2.Range("D17").Value = 0.47
$ws2.Range("E17").Value = 0.4

# --- Sheet1: label change ---
$ws1.Range("A7").Value = "All"

# --- Selections / active tab: Sheet1 ends up active with B16 selected,
#     correlation ends up inactive with E17 selected ---
$ws2.Range("E17").Select()
$ws1.Activate()
$ws1.Range("B16").Select()
